$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '62.749.31'
Set-TextValue "E2" '  -5.69%  '
Set-TextValue "D3" '3.102.91'
Set-TextValue "E3" '  -6.16%  '
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '558.49'
Set-TextValue "E5" '  -5.42%  '
Set-TextValue "D6" '161.26'
Set-TextValue "E6" '  -10.87%  '
Set-TextValue "E7" '  +0.04%  '
Set-TextValue "D8" '0.582'
Set-TextValue "E8" '  -8.67%  '
Set-TextValue "D9" '3.095.80'
Set-TextValue "E9" '  -6.31%  '
Set-TextValue "E10" '  -2.34%  '
Set-TextValue "D11" '0.115'
Set-TextValue "E11" '  -8.79%  '
Set-TextValue "D12" '0.377'
Set-TextValue "E12" '  -6.53%  '
Set-TextValue "D13" '3.636.36'
Set-TextValue "E13" '  -6.31%  '
Set-TextValue "E14" '  -2.13%  '
Set-TextValue "D15" '62.731.27'
Set-TextValue "E15" '  -5.73%  '
Set-TextValue "D16" '24.53'
Set-TextValue "E16" '  -8.17%  '
Set-TextValue "D17" '3.094.28'
Set-TextValue "E17" '  -6.80%  '
Set-TextValue "D18" '0.0000152'
Set-TextValue "E18" '  -7.25%  '
Set-TextValue "D19" '397.75'
Set-TextValue "D20" '12.31'
Set-TextValue "E20" '  -5.66%  '
Set-TextValue "D21" '5.11'
Set-TextValue "E21" '  -6.76%  '
Set-TextValue "D22" '7.03'
Set-TextValue "E22" '  -3.74%  '
Set-TextValue "E23" '  -0.18%  '
Set-TextValue "E24" '  -1.03%  '
Set-TextValue "D25" '67.44'
Set-TextValue "E25" '  -5.71%  '
Set-TextValue "D26" '0.195'
Set-TextValue "E26" '  -5.24%  '
Set-TextValue "D27" '0.477'
Set-TextValue "E27" '  -7.30%  '
Set-TextValue "D28" '0.0000100'
Set-TextValue "E28" '  -12.51%  '
Set-TextValue "B29" 'InternetComputer(DFINITY)'
Set-TextValue "C29" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D29" '8.59'
Set-TextValue "E29" '  -6.72%  '
Set-TextValue "B30" 'Binance-PegBSC-USD'
Set-TextValue "C30" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D30" '0.999'
Set-TextValue "E30" '  -0.01%  '
Set-TextValue "D32" '1.77'
Set-TextValue "E32" '  -7.93%  '
Set-TextValue "D33" '20.85'
Set-TextValue "E33" '  -6.81%  '
Set-TextValue "B34" 'Aptos'
Set-TextValue "C34" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D34" '6.19'
Set-TextValue "E34" '  -5.85%  '
Set-TextValue "B35" 'NEARProtocol'
Set-TextValue "C35" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D35" '4.77'
Set-TextValue "E35" '  -7.87%  '
Set-TextValue "D36" '153.07'
Set-TextValue "E36" '  -3.62%  '
Set-TextValue "D37" '1.09'
Set-TextValue "E37" '  -7.93%  '
Set-TextValue "E38" '  -8.36%  '
Set-TextValue "D39" '2.705.89'
Set-TextValue "E39" '  -5.66%  '
Set-TextValue "D40" '1.65'
Set-TextValue "E40" '  -8.23%  '
Set-TextValue "D41" '23.15'
Set-TextValue "E41" '  -12.20%  '
Set-TextValue "B42" 'Filecoin'
Set-TextValue "C42" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D42" '3.99'
Set-TextValue "E42" '  -7.75%  '
Set-TextValue "B43" 'OKB'
Set-TextValue "C43" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D43" '38.16'
Set-TextValue "E43" '  -4.08%  '
Set-TextValue "E44" '  -8.62%  '
Set-TextValue "D45" '0.0600'
Set-TextValue "E45" '  -8.36%  '
Set-TextValue "D46" '5.21'
Set-TextValue "E46" '  -12.15%  '
Set-TextValue "D47" '0.0253'
Set-TextValue "E47" '  -6.25%  '
Set-TextValue "B48" 'InjectiveProtocol'
Set-TextValue "C48" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D48" '20.68'
Set-TextValue "E48" '  -9.68%  '
Set-TextValue "B49" 'FirstDigitalUSD'
Set-TextValue "C49" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D49" '0.999'
Set-TextValue "E49" '  -0.05%  '
Set-TextValue "D50" '279.20'
Set-TextValue "E50" '  -10.90%  '
Set-TextValue "E51" '  -5.47%  '
